# hot fix for MS1 centroid Thermo .mzML files in Shotgun mode
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the H2O flag (column F) for the [M-H]-sn1 / [M-H]-sn2 /
# [M-H]-sn1-H2O / [M-H]-sn2-H2O rows.
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0

# Update the active selection to match the new state saved in the file.
$ws.Range("F7").Select()
